$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 20-26: mark the "Skipped" status entries as "Completed" now that
# the ROI Count workflow has processed them.
for ($i = 20; $i -le 26; $i++) {
    $ws.Cells.Item($i, 3).Value = "Completed"
}

# Append the newly processed PO numbers (ROI Count workflow) as rows 27-37.
$newRows = @(
    @("101865(0120860)", "SO-00022146", "Completed"),
    @("77218(0119760)",  "SO-00022147", "Completed"),
    @("421818(109260)",  "SO-00022151", "Completed"),
    @("60317(0125160)",  "SO-00022152", "Completed"),
    @("308312(0104460)", "SO-00022153", "Completed"),
    @("358566(0106160)", "SO-00022154", "Completed"),
    @("499485(0102060)", "SO-00022195", "Skipped"),
    @("281788(0104060)", "SO-00022196", "Skipped"),
    @("449581(0102860)", "SO-00022197", "Skipped"),
    @("215931(0115960)", "SO-00022206", "Skipped"),
    @("434129(0107560)", "SO-00022207", "Skipped")
)

$row = 27
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
